$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Objects_Toolbar")

# Row 3 ("Toolbar - Items") - switch locator type from xpath to id,
# object type from tab to movetoelement, and simplify the locator value
# to just the element id accordingly.
$ws.Range("C3").Value = "id"
$ws.Range("D3").Value = "movetoelement"
$ws.Range("E3").Value = "Ribbon.ListItem-title"
